$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting rows 6:70 down to 7:71
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new data record
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44530
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100106
$ws.Range("H6").Value = "Oleaginosos"
$ws.Range("I6").Value = 100106002
$ws.Range("J6").Value = "Palta"
$ws.Range("K6").Value = "Hass"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 28000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29000
$ws.Range("Q6").Value = "`$/bandeja 10 kilos"
$ws.Range("R6").Value = "Perú"
$ws.Range("S6").Value = 2900
$ws.Range("T6").Value = 10
